$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InterfaceTypes")

# InterfaceType command now requires a Unit: set the Unit (column H) for the
# two existing InterfaceType rows (ft1, ft2) to "tonnes".
$ws.Range("H2").Value = "tonnes"
$ws.Range("H3").Value = "tonnes"

# Leave the cursor where a user would naturally end up after filling H2:H3.
$ws.Range("H4").Select()
